$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 7835
$ws.Range("F4").Value = 7981
$ws.Range("F7").Value = 6912
$ws.Range("F8").Value = 3444
$ws.Range("F10").Value = 3758
$ws.Range("F11").Value = 54
$ws.Range("F12").Value = 61
$ws.Range("F13").Value = 56
$ws.Range("F15").Value = 95
$ws.Range("F17").Value = 13
$ws.Range("F19").Value = 335
$ws.Range("F23").Value = 3952
$ws.Range("F27").Value = 509
$ws.Range("F28").Value = 1561
$ws.Range("F29").Value = 86
$ws.Range("F31").Value = 2856
$ws.Range("F32").Value = 2036
$ws.Range("F37").Value = 3928
$ws.Range("F38").Value = 372
$ws.Range("F41").Value = 933
$ws.Range("F42").Value = 659
$ws.Range("F43").Value = 79
$ws.Range("F44").Value = 1511
$ws.Range("F45").Value = 254
$ws.Range("F48").Value = 672
$ws.Range("F49").Value = 16

$ws = $wb.Worksheets.Item(2)
$ws.Range("F10").Value = 27
$ws.Range("F15").Value = 1
$ws.Range("F17").Value = 419

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 149

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 149
$ws.Range("F6").Value = 7835
$ws.Range("F8").Value = 7981
$ws.Range("F11").Value = 6912
$ws.Range("F12").Value = 3444
$ws.Range("F14").Value = 3758
$ws.Range("F15").Value = 54
$ws.Range("F16").Value = 61
$ws.Range("F17").Value = 56
$ws.Range("F18").Value = 95
$ws.Range("F19").Value = 13
$ws.Range("F27").Value = 3952
$ws.Range("F31").Value = 509
$ws.Range("F32").Value = 1561
$ws.Range("F33").Value = 86
$ws.Range("F35").Value = 2856
$ws.Range("F36").Value = 2036
$ws.Range("F40").Value = 3928
$ws.Range("F41").Value = 372
$ws.Range("F45").Value = 1511
$ws.Range("F46").Value = 254
$ws.Range("F48").Value = 672
$ws.Range("F49").Value = 16
